$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Major Components ")

# The little 1.27*n reference table lived in E20:F26. Rows 22-26 still carry
# a non-default cell style (s="13") on their E/F cells, so clearing just the
# values/formulas there leaves behind empty-but-styled cells (and rows) --
# matching what Excel itself would keep on disk.
$ws.Range("E22:F26").ClearContents()

# Rows 20 and 21 never had any non-default formatting, so once their values
# are gone there is nothing left to preserve. Deleting them outright (and
# re-inserting two blank rows in their place) drops the rows from the saved
# XML entirely while keeping every row below at its original row number, so
# the sheet's used range stays A1:K26.
$ws.Rows("20:21").Delete()
$ws.Rows("20:21").Insert()

# Move the selection, as recorded by the author, onto G26:G27.
$ws.Range("G26:G27").Select()

$wb.Save()
